$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The task breakdown row for "T-08 / Debugging / 4 hours" in the first
# story block (SSDMS-11) is being removed; delete the entire row 15 and
# let Excel shift everything below it up, adjusting the SUM() formula,
# the shared G-column formula range, and the merged A/B story cells
# automatically.
$ws.Rows(15).Delete()

# The task IDs in column C are plain text labels (not formulas), so they
# don't renumber automatically when the row shifts up - relabel the
# now-shifted rows to keep the T-01..T-12 sequence contiguous.
$ws.Range("C7").Value  = "T-05"
$ws.Range("C11").Value = "T-06"
$ws.Range("C12").Value = "T-07"
$ws.Range("C13").Value = "T-08"
$ws.Range("C14").Value = "T-09"
$ws.Range("C15").Value = "T-10"
$ws.Range("C16").Value = "T-11"
$ws.Range("C17").Value = "T-12"
